# Auto update Excel log: append newly collected sensor readings
# to the PIR, Humidity, and Temperature sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(20, 1).Value = "'2026-01-28"
$ws.Cells.Item(20, 2).Value = "17:03:30"
$ws.Cells.Item(20, 3).Value = "17:00"
$ws.Cells.Item(20, 4).Value = "Bathroom"
$ws.Cells.Item(20, 5).Value = "No Motion"
$ws.Cells.Item(20, 6).Value = "Inactive"

$ws.Cells.Item(21, 1).Value = "'2026-01-28"
$ws.Cells.Item(21, 2).Value = "17:03:32"
$ws.Cells.Item(21, 3).Value = "17:00"
$ws.Cells.Item(21, 4).Value = "Bathroom"
$ws.Cells.Item(21, 5).Value = "No Motion"
$ws.Cells.Item(21, 6).Value = "Inactive"

$ws.Cells.Item(22, 1).Value = "'2026-01-28"
$ws.Cells.Item(22, 2).Value = "17:03:37"
$ws.Cells.Item(22, 3).Value = "17:00"
$ws.Cells.Item(22, 4).Value = "Bathroom"
$ws.Cells.Item(22, 5).Value = "No Motion"
$ws.Cells.Item(22, 6).Value = "Inactive"

$ws.Cells.Item(23, 1).Value = "'2026-01-28"
$ws.Cells.Item(23, 2).Value = "17:03:42"
$ws.Cells.Item(23, 3).Value = "17:00"
$ws.Cells.Item(23, 4).Value = "Bathroom"
$ws.Cells.Item(23, 5).Value = "No Motion"
$ws.Cells.Item(23, 6).Value = "Inactive"

$ws.Cells.Item(24, 1).Value = "'2026-01-28"
$ws.Cells.Item(24, 2).Value = "17:03:47"
$ws.Cells.Item(24, 3).Value = "17:00"
$ws.Cells.Item(24, 4).Value = "Bathroom"
$ws.Cells.Item(24, 5).Value = "No Motion"
$ws.Cells.Item(24, 6).Value = "Inactive"

$ws.Cells.Item(25, 1).Value = "'2026-01-28"
$ws.Cells.Item(25, 2).Value = "17:03:52"
$ws.Cells.Item(25, 3).Value = "17:00"
$ws.Cells.Item(25, 4).Value = "Bathroom"
$ws.Cells.Item(25, 5).Value = "No Motion"
$ws.Cells.Item(25, 6).Value = "Inactive"

$ws.Cells.Item(26, 1).Value = "'2026-01-28"
$ws.Cells.Item(26, 2).Value = "17:03:57"
$ws.Cells.Item(26, 3).Value = "17:00"
$ws.Cells.Item(26, 4).Value = "Bathroom"
$ws.Cells.Item(26, 5).Value = "No Motion"
$ws.Cells.Item(26, 6).Value = "Inactive"

$ws.Cells.Item(27, 1).Value = "'2026-01-28"
$ws.Cells.Item(27, 2).Value = "17:04:02"
$ws.Cells.Item(27, 3).Value = "17:00"
$ws.Cells.Item(27, 4).Value = "Bathroom"
$ws.Cells.Item(27, 5).Value = "No Motion"
$ws.Cells.Item(27, 6).Value = "Inactive"

$ws.Cells.Item(28, 1).Value = "'2026-01-28"
$ws.Cells.Item(28, 2).Value = "17:04:07"
$ws.Cells.Item(28, 3).Value = "17:00"
$ws.Cells.Item(28, 4).Value = "Bathroom"
$ws.Cells.Item(28, 5).Value = "No Motion"
$ws.Cells.Item(28, 6).Value = "Inactive"

$ws.Cells.Item(29, 1).Value = "'2026-01-28"
$ws.Cells.Item(29, 2).Value = "17:04:13"
$ws.Cells.Item(29, 3).Value = "17:00"
$ws.Cells.Item(29, 4).Value = "Bathroom"
$ws.Cells.Item(29, 5).Value = "No Motion"
$ws.Cells.Item(29, 6).Value = "Inactive"

$ws.Cells.Item(30, 1).Value = "'2026-01-28"
$ws.Cells.Item(30, 2).Value = "17:04:18"
$ws.Cells.Item(30, 3).Value = "17:00"
$ws.Cells.Item(30, 4).Value = "Bathroom"
$ws.Cells.Item(30, 5).Value = "No Motion"
$ws.Cells.Item(30, 6).Value = "Inactive"

$ws.Cells.Item(31, 1).Value = "'2026-01-28"
$ws.Cells.Item(31, 2).Value = "17:04:23"
$ws.Cells.Item(31, 3).Value = "17:00"
$ws.Cells.Item(31, 4).Value = "Bathroom"
$ws.Cells.Item(31, 5).Value = "No Motion"
$ws.Cells.Item(31, 6).Value = "Inactive"

$ws.Cells.Item(32, 1).Value = "'2026-01-28"
$ws.Cells.Item(32, 2).Value = "17:04:28"
$ws.Cells.Item(32, 3).Value = "17:00"
$ws.Cells.Item(32, 4).Value = "Bathroom"
$ws.Cells.Item(32, 5).Value = "No Motion"
$ws.Cells.Item(32, 6).Value = "Inactive"


$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(17, 1).Value = "'2026-01-28"
$ws.Cells.Item(17, 2).Value = "17:03:29"
$ws.Cells.Item(17, 3).Value = "17:00"
$ws.Cells.Item(17, 4).Value = "Bathroom"
$ws.Cells.Item(17, 5).Value = "'87.4%"
$ws.Cells.Item(17, 6).Value = "Active"

$ws.Cells.Item(18, 1).Value = "'2026-01-28"
$ws.Cells.Item(18, 2).Value = "17:03:29"
$ws.Cells.Item(18, 3).Value = "17:00"
$ws.Cells.Item(18, 4).Value = "Bathroom"
$ws.Cells.Item(18, 5).Value = "'87.4%"
$ws.Cells.Item(18, 6).Value = "Active"

$ws.Cells.Item(19, 1).Value = "'2026-01-28"
$ws.Cells.Item(19, 2).Value = "17:03:35"
$ws.Cells.Item(19, 3).Value = "17:00"
$ws.Cells.Item(19, 4).Value = "Bathroom"
$ws.Cells.Item(19, 5).Value = "'87.4%"
$ws.Cells.Item(19, 6).Value = "Active"

$ws.Cells.Item(20, 1).Value = "'2026-01-28"
$ws.Cells.Item(20, 2).Value = "17:03:39"
$ws.Cells.Item(20, 3).Value = "17:00"
$ws.Cells.Item(20, 4).Value = "Bathroom"
$ws.Cells.Item(20, 5).Value = "'86.5%"
$ws.Cells.Item(20, 6).Value = "Active"

$ws.Cells.Item(21, 1).Value = "'2026-01-28"
$ws.Cells.Item(21, 2).Value = "17:03:43"
$ws.Cells.Item(21, 3).Value = "17:00"
$ws.Cells.Item(21, 4).Value = "Bathroom"
$ws.Cells.Item(21, 5).Value = "'87.5%"
$ws.Cells.Item(21, 6).Value = "Active"

$ws.Cells.Item(22, 1).Value = "'2026-01-28"
$ws.Cells.Item(22, 2).Value = "17:03:47"
$ws.Cells.Item(22, 3).Value = "17:00"
$ws.Cells.Item(22, 4).Value = "Bathroom"
$ws.Cells.Item(22, 5).Value = "'87.5%"
$ws.Cells.Item(22, 6).Value = "Active"

$ws.Cells.Item(23, 1).Value = "'2026-01-28"
$ws.Cells.Item(23, 2).Value = "17:03:51"
$ws.Cells.Item(23, 3).Value = "17:00"
$ws.Cells.Item(23, 4).Value = "Bathroom"
$ws.Cells.Item(23, 5).Value = "'86.5%"
$ws.Cells.Item(23, 6).Value = "Active"

$ws.Cells.Item(24, 1).Value = "'2026-01-28"
$ws.Cells.Item(24, 2).Value = "17:03:59"
$ws.Cells.Item(24, 3).Value = "17:00"
$ws.Cells.Item(24, 4).Value = "Bathroom"
$ws.Cells.Item(24, 5).Value = "'86.5%"
$ws.Cells.Item(24, 6).Value = "Active"

$ws.Cells.Item(25, 1).Value = "'2026-01-28"
$ws.Cells.Item(25, 2).Value = "17:04:03"
$ws.Cells.Item(25, 3).Value = "17:00"
$ws.Cells.Item(25, 4).Value = "Bathroom"
$ws.Cells.Item(25, 5).Value = "'87.4%"
$ws.Cells.Item(25, 6).Value = "Active"

$ws.Cells.Item(26, 1).Value = "'2026-01-28"
$ws.Cells.Item(26, 2).Value = "17:04:07"
$ws.Cells.Item(26, 3).Value = "17:00"
$ws.Cells.Item(26, 4).Value = "Bathroom"
$ws.Cells.Item(26, 5).Value = "'87.5%"
$ws.Cells.Item(26, 6).Value = "Active"

$ws.Cells.Item(27, 1).Value = "'2026-01-28"
$ws.Cells.Item(27, 2).Value = "17:04:11"
$ws.Cells.Item(27, 3).Value = "17:00"
$ws.Cells.Item(27, 4).Value = "Bathroom"
$ws.Cells.Item(27, 5).Value = "'86.6%"
$ws.Cells.Item(27, 6).Value = "Active"

$ws.Cells.Item(28, 1).Value = "'2026-01-28"
$ws.Cells.Item(28, 2).Value = "17:04:15"
$ws.Cells.Item(28, 3).Value = "17:00"
$ws.Cells.Item(28, 4).Value = "Bathroom"
$ws.Cells.Item(28, 5).Value = "'87.5%"
$ws.Cells.Item(28, 6).Value = "Active"

$ws.Cells.Item(29, 1).Value = "'2026-01-28"
$ws.Cells.Item(29, 2).Value = "17:04:19"
$ws.Cells.Item(29, 3).Value = "17:00"
$ws.Cells.Item(29, 4).Value = "Bathroom"
$ws.Cells.Item(29, 5).Value = "'86.6%"
$ws.Cells.Item(29, 6).Value = "Active"

$ws.Cells.Item(30, 1).Value = "'2026-01-28"
$ws.Cells.Item(30, 2).Value = "17:04:23"
$ws.Cells.Item(30, 3).Value = "17:00"
$ws.Cells.Item(30, 4).Value = "Bathroom"
$ws.Cells.Item(30, 5).Value = "'87.5%"
$ws.Cells.Item(30, 6).Value = "Active"

$ws.Cells.Item(31, 1).Value = "'2026-01-28"
$ws.Cells.Item(31, 2).Value = "17:04:27"
$ws.Cells.Item(31, 3).Value = "17:00"
$ws.Cells.Item(31, 4).Value = "Bathroom"
$ws.Cells.Item(31, 5).Value = "'87.5%"
$ws.Cells.Item(31, 6).Value = "Active"


$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(17, 1).Value = "'2026-01-28"
$ws.Cells.Item(17, 2).Value = "17:03:29"
$ws.Cells.Item(17, 3).Value = "17:00"
$ws.Cells.Item(17, 4).Value = "Bathroom"
$ws.Cells.Item(17, 5).Value = "22.9C"
$ws.Cells.Item(17, 6).Value = "Active"

$ws.Cells.Item(18, 1).Value = "'2026-01-28"
$ws.Cells.Item(18, 2).Value = "17:03:29"
$ws.Cells.Item(18, 3).Value = "17:00"
$ws.Cells.Item(18, 4).Value = "Bathroom"
$ws.Cells.Item(18, 5).Value = "22.9C"
$ws.Cells.Item(18, 6).Value = "Active"

$ws.Cells.Item(19, 1).Value = "'2026-01-28"
$ws.Cells.Item(19, 2).Value = "17:03:35"
$ws.Cells.Item(19, 3).Value = "17:00"
$ws.Cells.Item(19, 4).Value = "Bathroom"
$ws.Cells.Item(19, 5).Value = "22.9C"
$ws.Cells.Item(19, 6).Value = "Active"

$ws.Cells.Item(20, 1).Value = "'2026-01-28"
$ws.Cells.Item(20, 2).Value = "17:03:39"
$ws.Cells.Item(20, 3).Value = "17:00"
$ws.Cells.Item(20, 4).Value = "Bathroom"
$ws.Cells.Item(20, 5).Value = "22.8C"
$ws.Cells.Item(20, 6).Value = "Active"

$ws.Cells.Item(21, 1).Value = "'2026-01-28"
$ws.Cells.Item(21, 2).Value = "17:03:43"
$ws.Cells.Item(21, 3).Value = "17:00"
$ws.Cells.Item(21, 4).Value = "Bathroom"
$ws.Cells.Item(21, 5).Value = "22.9C"
$ws.Cells.Item(21, 6).Value = "Active"

$ws.Cells.Item(22, 1).Value = "'2026-01-28"
$ws.Cells.Item(22, 2).Value = "17:03:47"
$ws.Cells.Item(22, 3).Value = "17:00"
$ws.Cells.Item(22, 4).Value = "Bathroom"
$ws.Cells.Item(22, 5).Value = "22.9C"
$ws.Cells.Item(22, 6).Value = "Active"

$ws.Cells.Item(23, 1).Value = "'2026-01-28"
$ws.Cells.Item(23, 2).Value = "17:03:51"
$ws.Cells.Item(23, 3).Value = "17:00"
$ws.Cells.Item(23, 4).Value = "Bathroom"
$ws.Cells.Item(23, 5).Value = "22.9C"
$ws.Cells.Item(23, 6).Value = "Active"

$ws.Cells.Item(24, 1).Value = "'2026-01-28"
$ws.Cells.Item(24, 2).Value = "17:03:59"
$ws.Cells.Item(24, 3).Value = "17:00"
$ws.Cells.Item(24, 4).Value = "Bathroom"
$ws.Cells.Item(24, 5).Value = "22.8C"
$ws.Cells.Item(24, 6).Value = "Active"

$ws.Cells.Item(25, 1).Value = "'2026-01-28"
$ws.Cells.Item(25, 2).Value = "17:04:03"
$ws.Cells.Item(25, 3).Value = "17:00"
$ws.Cells.Item(25, 4).Value = "Bathroom"
$ws.Cells.Item(25, 5).Value = "22.8C"
$ws.Cells.Item(25, 6).Value = "Active"

$ws.Cells.Item(26, 1).Value = "'2026-01-28"
$ws.Cells.Item(26, 2).Value = "17:04:07"
$ws.Cells.Item(26, 3).Value = "17:00"
$ws.Cells.Item(26, 4).Value = "Bathroom"
$ws.Cells.Item(26, 5).Value = "22.9C"
$ws.Cells.Item(26, 6).Value = "Active"

$ws.Cells.Item(27, 1).Value = "'2026-01-28"
$ws.Cells.Item(27, 2).Value = "17:04:11"
$ws.Cells.Item(27, 3).Value = "17:00"
$ws.Cells.Item(27, 4).Value = "Bathroom"
$ws.Cells.Item(27, 5).Value = "22.9C"
$ws.Cells.Item(27, 6).Value = "Active"

$ws.Cells.Item(28, 1).Value = "'2026-01-28"
$ws.Cells.Item(28, 2).Value = "17:04:15"
$ws.Cells.Item(28, 3).Value = "17:00"
$ws.Cells.Item(28, 4).Value = "Bathroom"
$ws.Cells.Item(28, 5).Value = "22.8C"
$ws.Cells.Item(28, 6).Value = "Active"

$ws.Cells.Item(29, 1).Value = "'2026-01-28"
$ws.Cells.Item(29, 2).Value = "17:04:19"
$ws.Cells.Item(29, 3).Value = "17:00"
$ws.Cells.Item(29, 4).Value = "Bathroom"
$ws.Cells.Item(29, 5).Value = "22.9C"
$ws.Cells.Item(29, 6).Value = "Active"

$ws.Cells.Item(30, 1).Value = "'2026-01-28"
$ws.Cells.Item(30, 2).Value = "17:04:23"
$ws.Cells.Item(30, 3).Value = "17:00"
$ws.Cells.Item(30, 4).Value = "Bathroom"
$ws.Cells.Item(30, 5).Value = "22.9C"
$ws.Cells.Item(30, 6).Value = "Active"

$ws.Cells.Item(31, 1).Value = "'2026-01-28"
$ws.Cells.Item(31, 2).Value = "17:04:27"
$ws.Cells.Item(31, 3).Value = "17:00"
$ws.Cells.Item(31, 4).Value = "Bathroom"
$ws.Cells.Item(31, 5).Value = "22.8C"
$ws.Cells.Item(31, 6).Value = "Active"

